# Auto-generated Excel COM-interop edit script
# Applies the IESO report refresh: updated CreatedAt timestamp and
# refreshed hourly intertie price values for columns V:Z (hours 20-24)
# across the affected data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CreatedAt timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-11-21T18:07:14"

# Row 4
$ws.Range("W4").Value = 131.21
$ws.Range("X4").Value = 401.58
$ws.Range("Y4").Value = 201.92
$ws.Range("Z4").Value = 154.6

# Row 5
$ws.Range("W5").Value = -57.66

# Row 6
$ws.Range("W6").Value = -14.24
$ws.Range("X6").Value = -11.86
$ws.Range("Y6").Value = -10.3
$ws.Range("Z6").Value = -7.27

# Row 8
$ws.Range("W8").Value = -6.13
$ws.Range("X8").Value = 203.99

# Row 9
$ws.Range("V9").Value = 130.84
$ws.Range("W9").Value = 130.49
$ws.Range("X9").Value = 406.16
$ws.Range("Y9").Value = 211.16
$ws.Range("Z9").Value = 163.5

# Row 10
$ws.Range("W10").Value = -57.66

# Row 11
$ws.Range("V11").Value = -13.61
$ws.Range("W11").Value = -14.96
$ws.Range("X11").Value = -7.28
$ws.Range("Y11").Value = -1.06
$ws.Range("Z11").Value = 1.64

# Row 13
$ws.Range("W13").Value = -6.13
$ws.Range("X13").Value = 203.99

# Row 14
$ws.Range("W14").Value = 188.15
$ws.Range("X14").Value = 406.16
$ws.Range("Y14").Value = 211.16
$ws.Range("Z14").Value = 163.5

# Row 16
$ws.Range("W16").Value = -14.96
$ws.Range("X16").Value = -7.28
$ws.Range("Y16").Value = -1.06
$ws.Range("Z16").Value = 1.64

# Row 18
$ws.Range("W18").Value = -6.13
$ws.Range("X18").Value = 203.99

# Row 19
$ws.Range("W19").Value = 116.33
$ws.Range("X19").Value = 400.1
$ws.Range("Y19").Value = 200.58
$ws.Range("Z19").Value = 153.72

# Row 20
$ws.Range("W20").Value = -70.39

# Row 21
$ws.Range("W21").Value = -16.39
$ws.Range("X21").Value = -13.34
$ws.Range("Y21").Value = -11.63
$ws.Range("Z21").Value = -8.15

# Row 23
$ws.Range("W23").Value = -6.13
$ws.Range("X23").Value = 203.98

# Row 24
$ws.Range("W24").Value = 129.05
$ws.Range("X24").Value = 400.1
$ws.Range("Y24").Value = 200.58
$ws.Range("Z24").Value = 153.72

# Row 25
$ws.Range("W25").Value = -57.66

# Row 26
$ws.Range("W26").Value = -16.39
$ws.Range("X26").Value = -13.34
$ws.Range("Y26").Value = -11.63
$ws.Range("Z26").Value = -8.15

# Row 28
$ws.Range("W28").Value = -6.13
$ws.Range("X28").Value = 203.98

# Row 29
$ws.Range("W29").Value = 114.22
$ws.Range("X29").Value = 398.28
$ws.Range("Y29").Value = 198.89
$ws.Range("Z29").Value = 152.56

# Row 30
$ws.Range("W30").Value = -70.39

# Row 31
$ws.Range("W31").Value = -18.5
$ws.Range("X31").Value = -15.16
$ws.Range("Y31").Value = -13.33
$ws.Range("Z31").Value = -9.31

# Row 33
$ws.Range("W33").Value = -6.13
$ws.Range("X33").Value = 203.99

# Row 34
$ws.Range("V34").Value = 131.32
$ws.Range("W34").Value = 190.71
$ws.Range("X34").Value = 412.19
$ws.Range("Y34").Value = 218.78
$ws.Range("Z34").Value = 170.93

# Row 36
$ws.Range("V36").Value = -13.13
$ws.Range("W36").Value = -12.4
$ws.Range("X36").Value = -1.25
$ws.Range("Y36").Value = 6.56
$ws.Range("Z36").Value = 9.06

# Row 38
$ws.Range("W38").Value = -6.13
$ws.Range("X38").Value = 203.98

# Row 39
$ws.Range("W39").Value = 131.21
$ws.Range("X39").Value = 401.58
$ws.Range("Y39").Value = 201.92
$ws.Range("Z39").Value = 154.6

# Row 40
$ws.Range("W40").Value = -57.66

# Row 41
$ws.Range("W41").Value = -14.24
$ws.Range("X41").Value = -11.86
$ws.Range("Y41").Value = -10.3
$ws.Range("Z41").Value = -7.27

# Row 43
$ws.Range("W43").Value = -6.13
$ws.Range("X43").Value = 203.99

# Row 44
$ws.Range("V44").Value = 144.16
$ws.Range("W44").Value = 202.27
$ws.Range("X44").Value = 412.19
$ws.Range("Y44").Value = 211.58
$ws.Range("Z44").Value = 162.36

# Row 46
$ws.Range("V46").Value = -0.29
$ws.Range("W46").Value = -0.83
$ws.Range("X46").Value = -1.25

# Row 48
$ws.Range("W48").Value = -6.13
$ws.Range("X48").Value = 203.98

# Row 49
$ws.Range("V49").Value = 154.82
$ws.Range("W49").Value = 217.89
$ws.Range("X49").Value = 426.1
$ws.Range("Y49").Value = 224.33
$ws.Range("Z49").Value = 171.65

# Row 51
$ws.Range("V51").Value = 10.37
$ws.Range("W51").Value = 14.79
$ws.Range("X51").Value = 12.66
$ws.Range("Y51").Value = 12.11
$ws.Range("Z51").Value = 9.779999999999999

# Row 53
$ws.Range("W53").Value = -6.13
$ws.Range("X53").Value = 203.99

# Row 54
$ws.Range("V54").Value = 142.74
$ws.Range("W54").Value = 204.37
$ws.Range("X54").Value = 419.03
$ws.Range("Y54").Value = 219.23
$ws.Range("Z54").Value = 167.91

# Row 56
$ws.Range("V56").Value = -1.71
$ws.Range("W56").Value = 1.26
$ws.Range("Y56").Value = 7.02
$ws.Range("Z56").Value = 6.04

# Row 58
$ws.Range("W58").Value = -6.13
$ws.Range("X58").Value = 203.99

# Row 59
$ws.Range("V59").Value = 150.78
$ws.Range("W59").Value = 211.83
$ws.Range("X59").Value = 421.03
$ws.Range("Y59").Value = 220.37
$ws.Range("Z59").Value = 168.97

# Row 61
$ws.Range("V61").Value = 6.33
$ws.Range("W61").Value = 8.720000000000001
$ws.Range("X61").Value = 7.6
$ws.Range("Y61").Value = 8.15
$ws.Range("Z61").Value = 7.1

# Row 63
$ws.Range("W63").Value = -6.13
$ws.Range("X63").Value = 203.99

# Row 64
$ws.Range("V64").Value = 153.67
$ws.Range("W64").Value = 215.99
$ws.Range("X64").Value = 424.92
$ws.Range("Y64").Value = 224.33
$ws.Range("Z64").Value = 171.65

# Row 66
$ws.Range("V66").Value = 9.220000000000001
$ws.Range("W66").Value = 12.88
$ws.Range("X66").Value = 11.49
$ws.Range("Y66").Value = 12.11
$ws.Range("Z66").Value = 9.779999999999999

# Row 68
$ws.Range("W68").Value = -6.13
$ws.Range("X68").Value = 203.98

# Row 69
$ws.Range("W69").Value = 218.37
$ws.Range("X69").Value = 426.81
$ws.Range("Y69").Value = 226.73
$ws.Range("Z69").Value = 173.68

# Row 71
$ws.Range("W71").Value = 15.27
$ws.Range("X71").Value = 13.37
$ws.Range("Y71").Value = 14.51
$ws.Range("Z71").Value = 11.81

# Row 73
$ws.Range("W73").Value = -6.13
$ws.Range("X73").Value = 203.99

# Row 74
$ws.Range("V74").Value = 151.26
$ws.Range("W74").Value = 213.2
$ws.Range("X74").Value = 421.48
$ws.Range("Y74").Value = 221.52
$ws.Range("Z74").Value = 169.5

# Row 76
$ws.Range("V76").Value = 6.81
$ws.Range("W76").Value = 10.09
$ws.Range("X76").Value = 8.050000000000001
$ws.Range("Y76").Value = 9.300000000000001
$ws.Range("Z76").Value = 7.63

# Row 78
$ws.Range("W78").Value = -6.13
$ws.Range("X78").Value = 203.99

# Row 79
$ws.Range("V79").Value = 152.05
$ws.Range("W79").Value = 214.37
$ws.Range("X79").Value = 422.78
$ws.Range("Y79").Value = 222.66
$ws.Range("Z79").Value = 170.26

# Row 81
$ws.Range("V81").Value = 7.59
$ws.Range("W81").Value = 11.27
$ws.Range("X81").Value = 9.34
$ws.Range("Y81").Value = 10.44

# Row 83
$ws.Range("W83").Value = -6.13
$ws.Range("X83").Value = 203.99

# Row 84
$ws.Range("V84").Value = 136.79
$ws.Range("W84").Value = 192.95
$ws.Range("X84").Value = 413.44
$ws.Range("Y84").Value = 213.71
$ws.Range("Z84").Value = 163.84

# Row 86
$ws.Range("V86").Value = -7.66
$ws.Range("W86").Value = -10.15
$ws.Range("X86").Value = 0
$ws.Range("Y86").Value = 1.5
$ws.Range("Z86").Value = 1.97

# Row 88
$ws.Range("W88").Value = -6.13
$ws.Range("X88").Value = 203.99

# Row 89
$ws.Range("V89").Value = 132.04
$ws.Range("W89").Value = 126.94
$ws.Range("X89").Value = 398.28
$ws.Range("Y89").Value = 198.89
$ws.Range("Z89").Value = 152.56

# Row 90
$ws.Range("W90").Value = -57.66

# Row 91
$ws.Range("V91").Value = -12.41
$ws.Range("W91").Value = -18.5
$ws.Range("X91").Value = -15.16
$ws.Range("Y91").Value = -13.33
$ws.Range("Z91").Value = -9.31

# Row 93
$ws.Range("W93").Value = -6.13
$ws.Range("X93").Value = 203.99
